$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 ("Speaker Introduction"): split the "Ft. Knox " run so that the
# trailing space after "Knox" is dropped and "Knox" becomes its own run.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2Content = $s2.Shapes.Item(2)
$s2Tr = $s2Content.TextFrame.TextRange
$s2Full = $s2Tr.Text
$s2Idx = $s2Full.IndexOf("Ft. Knox ")
$s2Sel = $s2Tr.Characters($s2Idx + 5, 5)
$s2Sel.Text = "Knox"

# ---------------------------------------------------------------------
# Slide 3 ("Motivation"): split "workflow implemented" into two runs and
# add a new paragraph placeholder for a workflow image.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3Content = $s3.Shapes.Item(2)
$s3Tr = $s3Content.TextFrame.TextRange
$s3Full = $s3Tr.Text
$s3Idx = $s3Full.IndexOf("implemented")
$s3Sel = $s3Tr.Characters($s3Idx + 1, 11)
$s3Sel.Text = "implemented"
$s3Content.TextFrame.TextRange.InsertAfter("`r[Workflow Image Here]")

# ---------------------------------------------------------------------
# Duplicate the DFA slide (slide 5) twice to become the new slide 6
# ("Event Driven Deterministic Finite Automata") and slide 7 ("Database
# Representation of an Event Driven DFA") BEFORE we rewrite slide 5's own
# text, so the duplicates still carry the original "DFA" wording/
# formatting context.
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)

$s6 = $s5.Duplicate().Item(1)
$s7 = $s6.Duplicate().Item(1)

# --- New slide 6: Event Driven Deterministic Finite Automata ---
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Event Driven Deterministic Finite Automata"
$s6Body = $s6.Shapes.Item(2).TextFrame.TextRange
$s6Body.Text = "A finite state machine that uses events instead of an alphabet.`rStates can convey meanings besides 'Accept' and 'Reject'`rCommonly used in embedded systems to track or implement complex hardware state."

# --- New slide 7: Database Representation of an Event Driven DFA ---
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Database Representation of an Event Driven DFA"
# Content placeholder on slide 7 stays blank, matching the duplicated
# slide 5 layout (no further edit required).

# ---------------------------------------------------------------------
# Slide 5 itself: rename the title to "Finite State Machine" (and let the
# title shrink back to its natural size now that the text is shorter),
# then fill in the body placeholder with descriptive text.
# ---------------------------------------------------------------------
$s5Title = $s5.Shapes.Item(1)
$s5Title.TextFrame.TextRange.Text = "Finite State Machine"
$s5Title.TextFrame2.AutoSize = 2

$s5Body = $s5.Shapes.Item(2).TextFrame.TextRange
$s5Body.Text = "A finite state machine is a mathematical abstraction that accepts or rejects a string based on states and transitions.`r[Insert FSM diagram that accepts a*b.]"
